$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Execution column (F) results
$ws.Range("F2").Value = "y"
$ws.Range("F3").Value = "y"
$ws.Range("F5").Value = "n"

# TC_05's TestDataSet changes from TD3 to TD1 (TD3 becomes unused and is dropped)
$ws.Range("D6").Value = "TD1"
$ws.Range("F6").Value = "y"

# Column E ("Browser") widened to fit its longest entry ("androidbrowser")
$ws.Columns.Item(5).ColumnWidth = 13.3

# Restore the active cell/selection to just F6
$ws.Range("F6").Select()
